# SD6501 Final Project Tests workbook - add "UAT" (User Acceptance Testing) sheet
# with feedback data collected from Melissa, Philip and Joseph, formatted and
# wrapped into an Excel table (Table3), mirroring the UnitTest/EspressoTests sheets.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the LAST tab ------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($sheetCount))
$ws.Name = "UAT"

# --- Column headers (entered left-to-right skipping A, then A itself, to mirror ---------
# --- the original authoring order that the shared-strings table was built in) -----------
$ws.Range("B1").Value = "Could you register for a new user account easily?"
$ws.Range("C1").Value = "Could you update your user login details?"
$ws.Range("D1").Value = "Was all the text readable to you?"
$ws.Range("E1").Value = "Could you create an Account?"
$ws.Range("F1").Value = "Could you create a Transaction?"
$ws.Range("G1").Value = "Was navigating the application intuitive? Did it not take too much thinking on what steps to do to accomplish your goals?"
$ws.Range("H1").Value = "Did the application do what you expected?"

# --- Melissa's feedback row -------------------------------------------------------------
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "Yes"
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "Yes"
$ws.Range("H2").Value = "Yes"
$ws.Range("G2").Value = "Yes browsing was intuitive and the application easy to use to accomplish goals like changing user details and performing transactions."

# --- Philip's feedback row ---------------------------------------------------------------
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Yes"
$ws.Range("H3").Value = "Yes"
$ws.Range("G3").Value = "Yes, easy to use"

# --- Joseph's feedback row ----------------------------------------------------------------
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "Yes"
$ws.Range("D4").Value = "Yes"
$ws.Range("E4").Value = "Yes"
$ws.Range("F4").Value = "Yes"
$ws.Range("H4").Value = "Yes"
$ws.Range("G4").Value = "It was easy and quick."

# --- "User" header + the three user names -------------------------------------------------
$ws.Range("A1").Value = "User"
$ws.Range("A2").Value = "Melissa"
$ws.Range("A3").Value = "Philip"
$ws.Range("A4").Value = "Joseph"

# --- Column widths, tuned to visually match the source workbook --------------------------
$ws.Columns.Item(2).ColumnWidth = 41.6   # Could you register...
$ws.Columns.Item(3).ColumnWidth = 35.3   # Could you update...
$ws.Columns.Item(4).ColumnWidth = 28.8   # Was all the text readable...
$ws.Columns.Item(5).ColumnWidth = 26.3   # Could you create an Account?
$ws.Columns.Item(6).ColumnWidth = 28.1   # Could you create a Transaction?
$ws.Columns.Item(7).ColumnWidth = 56.3   # Was navigating the application intuitive?...
$ws.Columns.Item(8).ColumnWidth = 36.6   # Did the application do what you expected?

# --- Header row formatting: centred, vertically centred, wrapped -------------------------
$ws.Range("A1:H1").HorizontalAlignment = -4108
$ws.Range("A1:H1").VerticalAlignment = -4108
$ws.Range("A1:H1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 43.2

# --- Data rows: top-aligned; the long answer from Melissa wraps onto multiple lines ------
$ws.Range("A2:H4").VerticalAlignment = -4160
$ws.Range("G2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 43.2

# --- Turn the range into an Excel Table (mirrors Table1/Table13 on the other sheets) -----
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:H4"), [System.Type]::Missing, 1)
$tbl.Name = "Table3"
$tbl.TableStyle = "TableStyleMedium2"

# --- Page setup to match the other worksheets ---------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- View state: scrolled right with F11 selected, and make UAT the active tab -----------
$ws.Range("F11").Select()
$ws.Activate()

Write-Host "UAT sheet added."
